# Daily attendance processing - re-sort the "Recorded By" (column G) list
# of names/emails alphabetically (case-insensitive) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 157
$col = 7  # column G = "Recorded By"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    $text = [string]$val

    if ($text -eq "") {
        continue
    }

    $parts = $text -split ",\s*"
    if ($parts.Count -le 1) {
        continue
    }

    $trimmed = $parts | ForEach-Object { $_.Trim() }
    $sorted = $trimmed | Sort-Object { $_.ToLower() }
    $joined = $sorted -join ", "

    if ($joined -ne $text) {
        $cell.Value = $joined
    }
}
